$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.157.34"
$ws.Range("E2").Value = "  +1.87%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.435.60"
$ws.Range("E3").Value = "  +2.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.52"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.72"
$ws.Range("E6").Value = "  -4.53%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").Value = "  +7.49%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.754"
$ws.Range("E9").Value = "  +11.19%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  +18.11%  "

# Row 11 - Avalanche
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.02"
$ws.Range("E11").Value = "  +0.61%  "

# Row 12 - TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.141"
$ws.Range("E12").Value = "  -0.40%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.85"
$ws.Range("E13").Value = "  +5.86%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.47"
$ws.Range("E14").Value = "  +4.21%  "

# Row 15 - now WrappedEther (was ShibaInu)
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.477.87"
$ws.Range("E15").Value = "  +3.17%  "

# Row 16 - now ShibaInu (was WrappedEther)
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000190"
$ws.Range("E16").Value = "  +48.95%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +3.20%  "

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.114.21"
$ws.Range("E18").Value = "  +1.82%  "

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.38"
$ws.Range("E19").Value = "  +3.51%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "390.83"
$ws.Range("E20").Value = "  +26.09%  "

# Row 21 - Litecoin
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "88.93"
$ws.Range("E21").Value = "  +5.58%  "

# Row 22 - ImmutableX
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.19"
$ws.Range("E22").Value = "  -1.10%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.37"
$ws.Range("E23").Value = "  +5.46%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  +2.89%  "

# Row 25 - EthereumClassic
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "32.11"
$ws.Range("E25").Value = "  +8.92%  "

# Row 26 - LEO
$ws.Range("E26").Value = "  +0.38%  "

# Row 27 - Filecoin
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.49"
$ws.Range("E27").Value = "  +1.44%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.69"
$ws.Range("E28").Value = "  +2.96%  "

# Row 29 - Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.73"
$ws.Range("E29").Value = "  +9.76%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "44.07"
$ws.Range("E30").Value = "  +6.78%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  -0.73%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  -0.43%  "

# Row 33 - Cosmos
$ws.Range("E33").Value = "  +4.37%  "

# Row 34 - Dai
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35 - VeChain
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0490"
$ws.Range("E35").Value = "  +1.89%  "

# Row 36 - OKB
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "52.31"
$ws.Range("E36").Value = "  +0.83%  "

# Row 37 - FirstDigitalUSD
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.09%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -1.67%  "

# Row 39 - now Stellar (was Stacks)
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.133"
$ws.Range("E39").Value = "  +7.75%  "

# Row 40 - now Stacks (was Stellar)
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  +0.86%  "

# Row 41 - TheGraph
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.311"
$ws.Range("E41").Value = "  +8.70%  "

# Row 42 - Monero
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.36"
$ws.Range("E42").Value = "  +3.18%  "

# Row 43 - ARBITRUM
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.97"
$ws.Range("E43").Value = "  -0.90%  "

# Row 44 - NEARProtocol
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.00"
$ws.Range("E44").Value = "  -0.42%  "

# Row 45 - Celestia
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.72"
$ws.Range("E45").Value = "  +0.36%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  +4.36%  "

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.96"
$ws.Range("E47").Value = "  +2.81%  "

# Row 48 - Maker
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.118.79"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49 - ApeXProtocol
$ws.Range("E49").Value = "  -0.10%  "

# Row 50 - ThetaToken
$ws.Range("E50").Value = "  +2.58%  "

# Row 51 - BEAM
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0366"
$ws.Range("E51").Value = "  +6.47%  "
